$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row data (columns B, C, E, G) ---
$data = @(
    @{ Row = 2;  B = "Marvo S6-280";           C = 2300;  E = "Marvo.jpg";     G = 2 },
    @{ Row = 3;  B = "DeepCool Tower Cooler";  C = 8500;  E = "D.jpg";         G = 2 },
    @{ Row = 4;  B = "Blue Shirt";             C = 500;   E = "shirt.jpg";     G = 6 },
    @{ Row = 5;  B = "Walton Indigo";          C = 4500;  E = "walton.jpg";    G = 4 },
    @{ Row = 6;  B = "Khichuri";               C = 200;   E = "khichuri.jpg";  G = 1 },
    @{ Row = 7;  B = "Castle";                 C = 12000; E = "castle.jpg";    G = 2 },
    @{ Row = 8;  B = "Gskill TridenZ";         C = 4500;  E = "gskill.jpg";    G = 23 },
    @{ Row = 9;  B = "AMD R5 S7K";             C = 32000; E = "AMDR5S7K.jpg";  G = 2 },
    @{ Row = 10; B = "Intel Core I5 13600K ";  C = 39000; E = "intel.jpg";     G = 2 },
    @{ Row = 11; B = "Holy Quran";             C = 3000;  E = "Quran.jpg";     G = 19 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 7).Value = $item.G
}

# --- Column widths (closest representable values; engine quantizes to 1/6-char px grid) ---
$ws.Columns.Item(2).ColumnWidth = 26.5
$ws.Columns.Item(3).ColumnWidth = 17.83

# --- Selection ---
$ws.Range("I15").Select()
